{"js": "// Update the worksheet date and all the division problems.\n// Each entry is [old text, new text]. Every old value is unique in the\n// document EXCEPT \"711\u00f76=\", which appears twice and maps to two\n// different replacements depending on which occurrence it is (handled\n// separately below, in document order).\nconst replacements = [\n  [\"2025-11-01 Saturday\", \"2025-11-02 Sunday\"],\n  [\"549\u00f79=\", \"810\u00f74=\"],\n  [\"260\u00f76=\", \"869\u00f79=\"],\n  [\"938\u00f77=\", \"152\u00f74=\"],\n  [\"678\u00f77=\", \"584\u00f77=\"],\n  [\"922\u00f75=\", \"480\u00f79=\"],\n  [\"718\u00f73=\", \"710\u00f75=\"],\n  [\"986\u00f75=\", \"168\u00f76=\"],\n  [\"262\u00f72=\", \"329\u00f73=\"],\n  [\"361\u00f72=\", \"304\u00f74=\"],\n  [\"552\u00f73=\", \"280\u00f72=\"],\n  [\"445\u00f78=\", \"205\u00f74=\"],\n  [\"604\u00f79=\", \"705\u00f77=\"],\n  [\"441\u00f76=\", \"223\u00f79=\"],\n  [\"448\u00f76=\", \"937\u00f73=\"],\n  [\"440\u00f79=\", \"188\u00f74=\"],\n  [\"263\u00f77=\", \"103\u00f74=\"],\n  [\"279\u00f74=\", \"315\u00f73=\"],\n  [\"369\u00f77=\", \"227\u00f78=\"],\n  [\"120\u00f77=\", \"451\u00f73=\"],\n  [\"811\u00f73=\", \"678\u00f72=\"],\n  [\"362\u00f72=\", \"733\u00f78=\"],\n  [\"401\u00f79=\", \"842\u00f79=\"],\n  [\"889\u00f73=\", \"268\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// \"711\u00f76=\" occurs twice: the first (document-order) occurrence becomes\n// \"318\u00f74=\" and the second becomes \"330\u00f75=\".\nconst dupResults = context.document.body.search(\"711\u00f76=\", { matchCase: true, matchWildcards: false });\ndupResults.load(\"items\");\nawait context.sync();\ndupResults.items[0].insertText(\"318\u00f74=\", Word.InsertLocation.replace);\ndupResults.items[1].insertText(\"330\u00f75=\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the worksheet date and all the division problems.\n$d = $word.ActiveDocument\n\n# --- Date line -------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"2025-11-01 Saturday\"\n$find.Replacement.Text = \"2025-11-02 Sunday\"\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n\n# --- Division problems -------------------------------------------------\n# The problems live in the single table on the page; only every 4th row\n# (1, 5, 9, 13, 17 in 1-based COM numbering) actually holds the 5 problem\n# cells, the rest are blank answer rows. Addressing by (row, column)\n# avoids ambiguity from the \"711\u00f76=\" value, which appears twice in the\n# sheet but maps to two different replacements.\n$newValues = @(\n    @(\"810\u00f74=\", \"869\u00f79=\", \"152\u00f74=\", \"584\u00f77=\", \"480\u00f79=\"),\n    @(\"710\u00f75=\", \"168\u00f76=\", \"329\u00f73=\", \"304\u00f74=\", \"280\u00f72=\"),\n    @(\"205\u00f74=\", \"705\u00f77=\", \"223\u00f79=\", \"937\u00f73=\", \"188\u00f74=\"),\n    @(\"103\u00f74=\", \"315\u00f73=\", \"227\u00f78=\", \"451\u00f73=\", \"318\u00f74=\"),\n    @(\"678\u00f72=\", \"733\u00f78=\", \"842\u00f79=\", \"268\u00f76=\", \"330\u00f75=\")\n)\n$rowNumbers = @(1, 5, 9, 13, 17)\n\n$tbl = $d.Tables.Item(1)\nfor ($g = 0; $g -lt $rowNumbers.Length; $g++) {\n    $rowNum = $rowNumbers[$g]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $tbl.Cell($rowNum, $c)\n        $cellRange = $cell.Range\n        # Trim the trailing end-of-cell marker so we only overwrite the\n        # visible text, which preserves the run/paragraph formatting.\n        $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)\n        $textRange.Text = $newValues[$g][$c - 1]\n    }\n}\n"}
